$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 434
$ws.Range("I2").Value = 1184
$ws.Range("J2").Value = 4983
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 1404
$ws.Range("M2").Value = 97
$ws.Range("N2").Value = 836
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 17
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 46
$ws.Range("S2").Value = 555
$ws.Range("T2").Value = 878
$ws.Range("U2").Value = 67
$ws.Range("V2").Value = 7730
$ws.Range("X2").Value = 7870
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 100
$ws.Range("AA2").Value = 60
